$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.688952333333333
$ws.Range("H2").Value = 14.066857
$ws.Range("I2").Value = 0.8812414385524715
$ws.Range("J2").Value = 0.8812414385524715
$ws.Range("M2").Value = 293.7577056666667
$ws.Range("N2").Value = 881.273117
$ws.Range("O2").Value = 0.9369756110667984
$ws.Range("P2").Value = 0.940660486426629
$ws.Range("Q2").Value = 1377.415879420363
$ws.Range("R2").Value = 12396.74291478327
$ws.Range("S2").Value = 0.8257017353850864
$ws.Range("T2").Value = 0.8289490002480702

$ws.Range("G3").Value = 4.688952333333333
$ws.Range("H3").Value = 14.066857
$ws.Range("I3").Value = 0.8812414385524715
$ws.Range("J3").Value = 0.8812414385524715
$ws.Range("O3").Value = 0.01525159481997056
$ws.Range("P3").Value = 0.01531157527761154
$ws.Range("Q3").Value = 22.42084921249367
$ws.Range("R3").Value = 201.787642912443
$ws.Range("S3").Value = 0.01344033735937028
$ws.Range("T3").Value = 0.01349319462414686

$ws.Range("G4").Value = 4.688952333333333
$ws.Range("H4").Value = 14.066857
$ws.Range("I4").Value = 0.8812414385524715
$ws.Range("J4").Value = 0.8812414385524715
$ws.Range("M4").Value = 7.906212666666666
$ws.Range("N4").Value = 23.718638
$ws.Range("O4").Value = 0.02521781829607561
$ws.Range("P4").Value = 0.02531699325449539
$ws.Range("Q4").Value = 37.07185433119621
$ws.Range("R4").Value = 333.6466889807659
$ws.Range("S4").Value = 0.02222298647238851
$ws.Range("T4").Value = 0.02231038355541473

$ws.Range("G5").Value = 4.688952333333333
$ws.Range("H5").Value = 14.066857
$ws.Range("I5").Value = 0.8812414385524715
$ws.Range("J5").Value = 0.8812414385524715
$ws.Range("M5").Value = 3.6844455
$ws.Range("N5").Value = 7.368891000000001
$ws.Range("O5").Value = 0.01175198303639443
$ws.Range("P5").Value = 0.007865466969060864
$ws.Range("Q5").Value = 17.2761893242645
$ws.Range("R5").Value = 103.657135945587
$ws.Range("S5").Value = 0.01035633443683647
$ws.Range("T5").Value = 0.006931375426702144

$ws.Range("G6").Value = 4.688952333333333
$ws.Range("H6").Value = 14.066857
$ws.Range("I6").Value = 0.8812414385524715
$ws.Range("J6").Value = 0.8812414385524715
$ws.Range("M6").Value = 3.386921
$ws.Range("N6").Value = 10.160763
$ws.Range("O6").Value = 0.01080299278076119
$ws.Range("P6").Value = 0.01084547807220323
$ws.Range("Q6").Value = 15.88111112576566
$ws.Range("R6").Value = 142.930000131891
$ws.Range("S6").Value = 0.009520044898789959
$ws.Range("T6").Value = 0.009557484698137662

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6318963333333333
$ws.Range("H7").Value = 1.895689
$ws.Range("I7").Value = 0.1187585614475285
$ws.Range("J7").Value = 0.1187585614475285
$ws.Range("M7").Value = 293.7577056666667
$ws.Range("N7").Value = 881.273117
$ws.Range("O7").Value = 0.9369756110667984
$ws.Range("P7").Value = 0.940660486426629
$ws.Range("Q7").Value = 185.6244170991792
$ws.Range("R7").Value = 1670.619753892613
$ws.Range("S7").Value = 0.1112738756817119
$ws.Range("T7").Value = 0.1117114861785589

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6318963333333333
$ws.Range("H8").Value = 1.895689
$ws.Range("I8").Value = 0.1187585614475285
$ws.Range("J8").Value = 0.1187585614475285
$ws.Range("O8").Value = 0.01525159481997056
$ws.Range("P8").Value = 0.01531157527761154
$ws.Range("Q8").Value = 3.021496360045667
$ws.Range("R8").Value = 27.193467240411
$ws.Range("S8").Value = 0.00181125746060028
$ws.Range("T8").Value = 0.001818380653464688

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6318963333333333
$ws.Range("H9").Value = 1.895689
$ws.Range("I9").Value = 0.1187585614475285
$ws.Range("J9").Value = 0.1187585614475285
$ws.Range("M9").Value = 7.906212666666666
$ws.Range("N9").Value = 23.718638
$ws.Range("O9").Value = 0.02521781829607561
$ws.Range("P9").Value = 0.02531699325449539
$ws.Range("Q9").Value = 4.995906794620221
$ws.Range("R9").Value = 44.963161151582
$ws.Range("S9").Value = 0.002994831823687103
$ws.Range("T9").Value = 0.003006609699080655

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6318963333333333
$ws.Range("H10").Value = 1.895689
$ws.Range("I10").Value = 0.1187585614475285
$ws.Range("J10").Value = 0.1187585614475285
$ws.Range("M10").Value = 3.6844455
$ws.Range("N10").Value = 7.368891000000001
$ws.Range("O10").Value = 0.01175198303639443
$ws.Range("P10").Value = 0.007865466969060864
$ws.Range("Q10").Value = 2.3281876018165
$ws.Range("R10").Value = 13.969125610899
$ws.Range("S10").Value = 0.001395648599557961
$ws.Range("T10").Value = 0.0009340915423587203

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.6318963333333333
$ws.Range("H11").Value = 1.895689
$ws.Range("I11").Value = 0.1187585614475285
$ws.Range("J11").Value = 0.1187585614475285
$ws.Range("M11").Value = 3.386921
$ws.Range("N11").Value = 10.160763
$ws.Range("O11").Value = 0.01080299278076119
$ws.Range("P11").Value = 0.01084547807220323
$ws.Range("Q11").Value = 2.140182961189666
$ws.Range("R11").Value = 19.261646650707
$ws.Range("S11").Value = 0.001282947881971235
$ws.Range("T11").Value = 0.00128799337406557
